$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell address -> new price value.
# Prices are stored as text in the sheet (t="inlineStr"), so we force
# a text number format before assignment to stop Excel auto-converting
# the numeric-looking string into a real number, then restore the
# default "Normal" style so no extra formatting is left behind.
$updates = @{
    "D2" = "254.34"
    "D3" = "23.79"
    "D4" = "6.119"
    "D5" = "0.06003"
    "D6" = "3.440"
    "D7" = "6.605"
    "D8" = "1.326"
    "D9" = "0.8025"
    "D10" = "0.1521"
    "D11" = "0.07963"
    "D12" = "0.03341"
    "D13" = "0.03110"
    "D14" = "0.09297"
    "D15" = "3.600"
    "D16" = "0.001684"
    "D17" = "0.04791"
    "D18" = "0.0006117"
    "D19" = "0.006274"
    "D20" = "0.005893"
    "D21" = "0.001089"
    "D22" = "0.0001508"
    "D23" = "3.664"
    "D26" = "0.1265"
    "D27" = "0.0006515"
    "D40" = "0.04475"
    "D41" = "0.007088"
    "D42" = "0.1072"
    "D43" = "0.003367"
    "D44" = "0.01054"
    "D45" = "0.002475"
    "D46" = "0.00005929"
    "D47" = "0.00000000754"
    "D48" = "0.7044"
    "D49" = "0.1210"
    "D50" = "0.00002111"
    "D51" = "0.01015"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
